$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '27.701.86'
Set-TextValue 'E2' '  -0.39%  '

# Row 3
Set-TextValue 'D3' '1.590.01'
Set-TextValue 'E3' '  -2.22%  '

# Row 4
Set-TextValue 'E4' '  +0.69%  '

# Row 5
Set-TextValue 'D5' '207.29'
Set-TextValue 'E5' '  -1.68%  '

# Row 6
Set-TextValue 'E6' '  -2.81%  '

# Row 7
Set-TextValue 'E7' '  +0.68%  '

# Row 8
Set-TextValue 'D8' '22.23'
Set-TextValue 'E8' '  -3.95%  '

# Row 9
Set-TextValue 'E9' '  -1.57%  '

# Row 11
Set-TextValue 'D11' '0.0868'
Set-TextValue 'E11' '  -1.11%  '

# Row 12
Set-TextValue 'D12' '1.816.07'
Set-TextValue 'E12' '  -2.20%  '

# Row 13
Set-TextValue 'D13' '1.582.32'
Set-TextValue 'E13' '  -3.10%  '

# Row 14
Set-TextValue 'D14' '3.86'
Set-TextValue 'E14' '  -3.56%  '

# Row 15
Set-TextValue 'D15' '0.530'
Set-TextValue 'E15' '  -4.42%  '

# Row 16
Set-TextValue 'B16' 'WrappedBTC'
Set-TextValue 'C16' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D16' '27.694.86'
Set-TextValue 'E16' '  -0.48%  '

# Row 17
Set-TextValue 'B17' 'Litecoin'
Set-TextValue 'C17' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D17' '63.51'
Set-TextValue 'E17' '  -2.08%  '

# Row 18
Set-TextValue 'D18' '220.16'
Set-TextValue 'E18' '  -3.44%  '

# Row 19
Set-TextValue 'E19' '  -2.98%  '

# Row 20
Set-TextValue 'E20' '  -3.56%  '

# Row 22
Set-TextValue 'E22' '  -4.52%  '

# Row 23
Set-TextValue 'D23' '9.59'
Set-TextValue 'E23' '  -3.24%  '

# Row 24
Set-TextValue 'D24' '1.98'
Set-TextValue 'E24' '  -3.59%  '

# Row 25
Set-TextValue 'D25' '153.19'
Set-TextValue 'E25' '  -1.34%  '

# Row 26
Set-TextValue 'D26' '6.87'
Set-TextValue 'E26' '  -0.80%  '

# Row 27
Set-TextValue 'E27' '  +0.68%  '

# Row 28
Set-TextValue 'D28' '15.17'
Set-TextValue 'E28' '  -1.77%  '

# Row 29
Set-TextValue 'E29' '  -4.11%  '

# Row 30
Set-TextValue 'D30' '1.16'
Set-TextValue 'E30' '  -1.62%  '

# Row 31
Set-TextValue 'D31' '0.0469'
Set-TextValue 'E31' '  -2.09%  '

# Row 32
Set-TextValue 'E32' '  -5.11%  '

# Row 33
Set-TextValue 'D33' '1.373.00'
Set-TextValue 'E33' '  -2.50%  '

# Row 34
Set-TextValue 'E34' '  -5.40%  '

# Row 35
Set-TextValue 'E35' '  -4.16%  '

# Row 36
Set-TextValue 'E36' '  -1.78%  '

# Row 37
Set-TextValue 'E37' '  -0.65%  '

# Row 38
Set-TextValue 'E38' '  -0.95%  '

# Row 39
Set-TextValue 'D39' '0.539'
Set-TextValue 'E39' '  -2.43%  '

# Row 40
Set-TextValue 'D40' '0.825'
Set-TextValue 'E40' '  -2.75%  '

# Row 41
Set-TextValue 'E41' '  +0.56%  '

# Row 42
Set-TextValue 'D42' '0.970'
Set-TextValue 'E42' '  -3.05%  '

# Row 43
Set-TextValue 'D43' '64.44'
Set-TextValue 'E43' '  -1.87%  '

# Row 44
Set-TextValue 'E44' '  +2.71%  '

# Row 45
Set-TextValue 'E45' '  -3.02%  '

# Row 46
Set-TextValue 'E46' '  -4.08%  '

# Row 47
Set-TextValue 'D47' '1.726.87'
Set-TextValue 'E47' '  -2.23%  '

# Row 48
Set-TextValue 'D48' '87.89'
Set-TextValue 'E48' '  -0.67%  '

# Row 49
Set-TextValue 'E49' '  +11.79%  '

# Row 50
Set-TextValue 'D50' '0.0971'
Set-TextValue 'E50' '  -3.75%  '

# Row 51
Set-TextValue 'E51' '  -1.02%  '
